$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Individual cell value updates
$ws.Range("N5").Value = 9
$ws.Range("N16").Value = 17
$ws.Range("J18").Value = 2.37
$ws.Range("K21").Value = 2.37

# Row 22 updates
$ws.Range("G22").Value = 1.73
$ws.Range("I22").Value = 4.75
$ws.Range("J22").Value = 2.4
$ws.Range("L22").Value = 5.5
$ws.Range("N22").Value = 9
$ws.Range("O22").Value = 1.36
$ws.Range("P22").Value = 3
$ws.Range("W22").Value = 6
$ws.Range("X22").Value = 7.5
$ws.Range("Y22").Value = 8.5
$ws.Range("Z22").Value = 13
$ws.Range("AC22").Value = 9
$ws.Range("AE22").Value = 19
$ws.Range("AF22").Value = 67
$ws.Range("AI22").Value = 23
$ws.Range("AN22").Value = 3.6
$ws.Range("AS22").Value = 201
$ws.Range("AU22").Value = 9
$ws.Range("AW22").Value = 6.5
$ws.Range("AX22").Value = 29

# New row 25 data
$ws.Range("A25").Value = "QH0N67zQ"
$ws.Range("B25").Value = "16/11/2024"
$ws.Range("C25").Value = "20:00"
$ws.Range("D25").Value = "VENEZUELA - LIGA FUTVE"
$ws.Range("E25").Value = "Dep. Tachira"
$ws.Range("F25").Value = "Zamora"
$ws.Range("G25").Value = 1.3
$ws.Range("H25").Value = 4.3
$ws.Range("I25").Value = 11
$ws.Range("J25").Value = 1.8
$ws.Range("K25").Value = 2.22
$ws.Range("L25").Value = 8.75
$ws.Range("M25").Value = 1.01
$ws.Range("N25").Value = 7.6
$ws.Range("O25").Value = 1.26
$ws.Range("P25").Value = 3.15
$ws.Range("Q25").Value = 1.83
$ws.Range("R25").Value = 1.87
$ws.Range("S25").Value = 1.4
$ws.Range("T25").Value = 2.52
$ws.Range("U25").Value = 2.18
$ws.Range("V25").Value = 1.53
$ws.Range("W25").Value = 5.6
$ws.Range("X25").Value = 5.4
$ws.Range("Y25").Value = 8.75
$ws.Range("Z25").Value = 7.6
$ws.Range("AA25").Value = 12
$ws.Range("AB25").Value = 35
$ws.Range("AC25").Value = 10
$ws.Range("AD25").Value = 9
$ws.Range("AE25").Value = 24
$ws.Range("AF25").Value = 150
$ws.Range("AG25").Value = 201
$ws.Range("AH25").Value = 26
$ws.Range("AI25").Value = 90
$ws.Range("AJ25").Value = 35
$ws.Range("AK25").Value = 450
$ws.Range("AL25").Value = 175
$ws.Range("AM25").Value = 120
$ws.Range("AN25").Value = 2.92
$ws.Range("AO25").Value = 5.9
$ws.Range("AP25").Value = 18.5
$ws.Range("AQ25").Value = 16.5
$ws.Range("AR25").Value = 55
$ws.Range("AS25").Value = 300
$ws.Range("AT25").Value = 2.55
$ws.Range("AU25").Value = 9
$ws.Range("AV25").Value = 100
$ws.Range("AW25").Value = 10.5
$ws.Range("AX25").Value = 65
$ws.Range("AY25").Value = 60
$ws.Range("AZ25").Value = 600
$ws.Range("BA25").Value = 500
$ws.Range("BB25").Value = 51
$ws.Range("BC25").Value = 51
$ws.Range("BD25").Value = 51
